$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "B" and "C" quarter rows within each year group (the data rows
# were reordered while the row labels A/B/C/D stayed anchored to the same
# row numbers, so row 3 now carries what used to be row 4's values, etc.)
$swapPairs = @(
    @(3, 4),
    @(7, 8),
    @(11, 12),
    @(15, 16)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    # NOTE: use .Value2 (not .Value) to read a real scalar out of the cell -
    # in this host .Value round-trips a live property-accessor instead of a
    # snapshot, so re-assigning it later re-reads the (by-then-mutated) cell.
    $label1 = $ws.Range("A$r1").Value2
    $b1 = $ws.Range("B$r1").Value2
    $c1 = $ws.Range("C$r1").Value2
    $d1 = $ws.Range("D$r1").Value2
    $e1 = $ws.Range("E$r1").Value2

    $label2 = $ws.Range("A$r2").Value2
    $b2 = $ws.Range("B$r2").Value2
    $c2 = $ws.Range("C$r2").Value2
    $d2 = $ws.Range("D$r2").Value2
    $e2 = $ws.Range("E$r2").Value2

    $ws.Range("A$r1").Value2 = $label2
    $ws.Range("B$r1").Value2 = $b2
    $ws.Range("C$r1").Value2 = $c2
    $ws.Range("D$r1").Value2 = $d2
    $ws.Range("E$r1").Value2 = $e2

    $ws.Range("A$r2").Value2 = $label1
    $ws.Range("B$r2").Value2 = $b1
    $ws.Range("C$r2").Value2 = $c1
    $ws.Range("D$r2").Value2 = $d1
    $ws.Range("E$r2").Value2 = $e1
}

# Remove the now-redundant "挖掘机产销率" (F) and "挖掘机销售量" (G) columns
# entirely - shifts the used range back down to A1:E17.
$ws.Range("F1:G17").Delete()
